# Fruta / hortaliza, semanal
# Insert one new weekly data row above row 188 (pushing existing rows down),
# matching the new "Pehuenche" / "1a nueva(o)" observation recorded for
# 2021-11-23 at Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 188, shifting rows 188:278 down to 189:279
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row with the new record
$ws.Cells.Item(188, 1).Value = 4
$ws.Cells.Item(188, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(188, 3).Value = "Los Lagos"
$ws.Cells.Item(188, 4).Value = 44523
$ws.Cells.Item(188, 5).Value = 10
$ws.Cells.Item(188, 6).Value = 100114001
$ws.Cells.Item(188, 7).Value = "Papa"
$ws.Cells.Item(188, 8).Value = "Pehuenche"
$ws.Cells.Item(188, 9).Value = "1a nueva(o)"
$ws.Cells.Item(188, 10).Value = 600
$ws.Cells.Item(188, 11).Value = 13000
$ws.Cells.Item(188, 12).Value = 14000
$ws.Cells.Item(188, 13).Value = 13500
$ws.Cells.Item(188, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(188, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(188, 16).Value = 540
$ws.Cells.Item(188, 17).Value = 25
$ws.Cells.Item(188, 18).Value = "Hortaliza"
